$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells with the same style as the existing header row (s="1")
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the team record (Wins/Losses/Ties) for every data row
for ($r = 2; $r -le 44; $r++) {
    $ws.Cells.Item($r, 30).Value = 86
    $ws.Cells.Item($r, 31).Value = 76
    $ws.Cells.Item($r, 32).Value = 0
}
